# Generate Report for Archive
#
# The localization status for the "Test`1.md" file moved from
# "Ready for handoff" to "In Translation" for the de-de locale, and a
# Handoff name/timestamp was recorded for the new handoff round.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is Test`1.md, column E is the de-de status.
$overview.Range("E3").Value = "In Translation"

# de-de sheet: row 3 is Test`1.md, column C is Status, column I is
# "Lastest Handoff Name".
$dede.Range("C3").Value = "In Translation"
$dede.Range("I3").Value = "LocaleLowerCaseTest_HT_OL#Test1#20171103T035934"

# Widen the "Lastest Handoff Name" column to fit the new value (matches
# the width already used for the same column on the zh-cn sheet).
$dede.Columns.Item(9).ColumnWidth = 39.17
